$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$data = @(
    @(2, 0.0065692528824474366, 0.0049064732260763215),
    @(3, 9.6843820391987609, 0.008297765709917861),
    @(4, 0.046013031153623644, 0.0086062637502700174),
    @(5, 0.0037520512462482216, 0.0037636018467475272),
    @(6, 0.0036240228190348679, 0.0035137938526339481),
    @(7, 0.0079995714533636977, 0.0074709129920977604),
    @(8, 0.0077970321493506904, 0.0077851203517530663),
    @(9, 0.010570784760577111, 0.0072328182928707617),
    @(10, 0.010620402990019288, 0.0071167893531227638),
    @(11, 0.0062645743086167651, 0.0047727097477873398),
    @(12, 0.004609653326248824, 0.0045762281061653611),
    @(13, 0.0032006495552294287, 0.0031888921719043814),
    @(14, 0.0078953879747093794, 0.0071178594994428697),
    @(15, 0.0065924123641591488, 0.0053631834056428879),
    @(16, 0.0079799931157258676, 0.0080488427912783315),
    @(17, 0.0096373805275684413, 0.009427085020537072),
    @(18, 0.010052511632793446, 0.010017308286435137),
    @(19, 0.0051983574061072934, 0.0045062152542311639),
    @(20, 0.0050694161133063185, 0.0049986629586597255),
    @(21, 0.0064283343288796709, 0.0048188426224936255),
    @(22, 3.2923100975761921, 0.0085042561367234745),
    @(23, 0.05516526145208759, 0.010405576398268022),
    @(24, 0.0035176534716765194, 0.0035276831052710902),
    @(25, 0.0035003939336604806, 0.0034060262998012705),
    @(26, 0.0077042048232094443, 0.0072198938249473351),
    @(27, 0.0075898559418690007, 0.0075755584032157416),
    @(28, 0.010484476510234497, 0.0071722958651268629),
    @(29, 0.010353487341055005, 0.0069448151168824386),
    @(30, 0.0060586258448357133, 0.004618363243502354),
    @(31, 0.0044039700123535939, 0.0043701686391642873),
    @(32, 0.0030341255771608004, 0.0030212114399744395),
    @(33, 0.0074866033143364003, 0.0067869438716258073),
    @(34, 0.0063933108080743255, 0.0052077957760795776),
    @(35, 0.0078126337919887951, 0.0078750084790124852),
    @(36, 0.0094755177596258042, 0.0092678786726013134),
    @(37, 0.0099294611084738121, 0.0098883058884528148),
    @(38, 0.0051703265345527339, 0.0044811024169603784),
    @(39, 0.0050030156051434106, 0.0049303694478308674),
    @(40, 0.0063406158384813199, 0.0047606878781183033),
    @(41, 2.658261677986562, 0.0085030866946273834),
    @(42, 0.05457549115409966, 0.010296072882860824),
    @(43, 0.0035308084891429274, 0.0035407690733959625),
    @(44, 0.003488023559308009, 0.0033955045816051826),
    @(45, 0.0076366032547371737, 0.0071606637830319286),
    @(46, 0.0075286467339499728, 0.0075145030454718525),
    @(47, 0.010413832356695274, 0.0071242633817524036),
    @(48, 0.010248078739797401, 0.0068759109575553187),
    @(49, 0.0060129848272251069, 0.0045840927010338929),
    @(50, 0.0043468734698584813, 0.0043134718114194448),
    @(51, 0.002971570819486236, 0.0029588292901527682),
    @(52, 0.0074307548053437504, 0.0067410265893634667),
    @(53, 0.0063583907038028825, 0.0051803217287598787),
    @(54, 0.0077723177215585811, 0.0078340570322108834),
    @(55, 0.0094544477966035465, 0.009247500386868758),
    @(56, 0.0098928826349062686, 0.0098514765817128043),
    @(57, 0.0051799367762708527, 0.00448942247824974),
    @(58, 0.0049623963364967012, 0.0048902045094949918),
    @(59, 0.0062690726253979883, 0.0047121465940318987),
    @(60, 2.2887712582606863, 0.0084886310664350588),
    @(61, 0.054065035672580801, 0.010200962466271384),
    @(62, 0.0035171567343556362, 0.0035270221798835324),
    @(63, 0.0034853096540172802, 0.0033936918632493799),
    @(64, 0.0075786431365619375, 0.0071090828525772658),
    @(65, 0.0074774520694248095, 0.0074634160132823731),
    @(66, 0.010350362025417325, 0.0070810181026234962),
    @(67, 0.010158573595590368, 0.0068170646511555022),
    @(68, 0.0059674951760433428, 0.0045497980997743397),
    @(69, 0.0043022483886119784, 0.0042691589337456931),
    @(70, 0.0029273394350090524, 0.0029147220050261902),
    @(71, 0.0073802051543054329, 0.0066983589773226774),
    @(72, 0.0063286445317991581, 0.0051566584189946243),
    @(73, 0.0077356786335194992, 0.0077969154472332528),
    @(74, 0.0094311148834714392, 0.0092248051092712682),
    @(75, 0.0098586958123180624, 0.0098171607966421881),
    @(76, 0.005193118958946712, 0.0045008048153158094),
    @(77, 0.0049281149094526086, 0.0048563284859757094)
)
foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
}
